$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mention text in G3
$ws.Range("G3").Value = " @xauxauxa"

# Insert a new "#! END_ROW" marker cell in I2, copying the formatting
# from the existing I1 marker cell
$ws.Range("I1").Copy($ws.Range("I2"))

# Move the active selection to I4
$ws.Range("I4").Select()
